$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (string) cell type for the columns being updated, so that values
# like percentages and plain numbers are stored as literal text, matching the
# original inline-string cell contents instead of being auto-converted by Excel.

$rngD = $ws.Range("D2:D50")
$rngD.NumberFormat = "@"

$rngE = $ws.Range("E2:E48")
$rngE.NumberFormat = "@"

$rngG = $ws.Range("G2:G51")
$rngG.NumberFormat = "@"

$ws.Range("D2").Value = "302.43"
$ws.Range("D3").Value = "32.17"
$ws.Range("D4").Value = "5.011"
$ws.Range("D5").Value = "0.07912"
$ws.Range("D6").Value = "2.117"
$ws.Range("D7").Value = "7.863"
$ws.Range("D8").Value = "3.811"
$ws.Range("D9").Value = "0.9292"
$ws.Range("D10").Value = "0.1758"
$ws.Range("D11").Value = "0.07921"
$ws.Range("D12").Value = "0.08774"
$ws.Range("D13").Value = "0.03130"
$ws.Range("D14").Value = "0.1005"
$ws.Range("D15").Value = "0.001539"
$ws.Range("D16").Value = "0.005954"
$ws.Range("D17").Value = "3.470"
$ws.Range("D20").Value = "0.1289"
$ws.Range("D21").Value = "4.157"
$ws.Range("D23").Value = "0.04593"
$ws.Range("D25").Value = "0.004465"
$ws.Range("D39").Value = "0.01725"
$ws.Range("D40").Value = "0.04818"
$ws.Range("D41").Value = "0.007512"
$ws.Range("D42").Value = "0.1367"
$ws.Range("D43").Value = "0.002364"
$ws.Range("D44").Value = "0.01020"
$ws.Range("D45").Value = "0.00006012"
$ws.Range("D47").Value = "0.003395"
$ws.Range("D49").Value = "0.00002103"
$ws.Range("D50").Value = "0.0002003"

$ws.Range("E2").Value = "0.78%"
$ws.Range("E3").Value = "1.18%"
$ws.Range("E4").Value = "-2.62%"
$ws.Range("E5").Value = "-2.85%"
$ws.Range("E6").Value = "-13.79%"
$ws.Range("E7").Value = "0.30%"
$ws.Range("E8").Value = "-1.61%"
$ws.Range("E9").Value = "0.38%"
$ws.Range("E10").Value = "-0.11%"
$ws.Range("E11").Value = "7.11%"
$ws.Range("E12").Value = "-2.82%"
$ws.Range("E13").Value = "3.30%"
$ws.Range("E14").Value = "0.38%"
$ws.Range("E15").Value = "1.27%"
$ws.Range("E16").Value = "-0.33%"
$ws.Range("E17").Value = "-3.37%"
$ws.Range("E18").Value = "-0.34%"
$ws.Range("E19").Value = "0.73%"
$ws.Range("E20").Value = "-3.77%"
$ws.Range("E21").Value = "1.63%"
$ws.Range("E22").Value = "6.70%"
$ws.Range("E23").Value = "-0.82%"
$ws.Range("E24").Value = "-0.73%"
$ws.Range("E25").Value = "-1.78%"
$ws.Range("E26").Value = "4.29%"
$ws.Range("E39").Value = "-2.13%"
$ws.Range("E41").Value = "8.69%"
$ws.Range("E42").Value = "-0.95%"
$ws.Range("E43").Value = "6.93%"
$ws.Range("E44").Value = "3.40%"
$ws.Range("E45").Value = "-4.71%"
$ws.Range("E48").Value = "2.24%"

# Every row from 2-51 gets "Hora" bumped from "20" to "21".
$rngG.Value = "21"

$rngD.Style = "Normal"
$rngE.Style = "Normal"
$rngG.Style = "Normal"
